$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# feat: new features, fix css and link images
# ---------------------------------------------------------------
# 1) The project screenshots were re-hosted from imgur.com to
#    postimg.cc - update the "imagem" (image link) column values.
$ws.Cells.Replace("https://i.imgur.com/0G6m8wt.png", "https://i.postimg.cc/kGzWSHLb/scrapy.png")
$ws.Cells.Replace("https://i.imgur.com/0LnJx5e.png", "https://i.postimg.cc/k5qKWKYf/email.png")
$ws.Cells.Replace("https://i.imgur.com/y8PIg61.png", "https://i.postimg.cc/YCGp2C6R/XML.png")

# 2) fix css: the header cells E1 ("link") and F1 ("imagem") were
#    using a plain/unbordered style while the rest of the header
#    row (A1:D1) uses a bordered style. Copy D1's formatting onto
#    E1:F1 so the whole header row is visually consistent.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
